# "Generate Report for Handback"
# Refresh the handoff/handback timestamps recorded for the 291ea829... item
# in both the zh-cn and de-de report sheets.

$wb = $excel.ActiveWorkbook

# zh-cn sheet: Correspond Handoff Datetime (E2) / Correspond Handback DateTime (H2)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-17 18:54:27"
$wsZhCn.Range("H2").Value = "2016-03-17 18:54:47"

# de-de sheet: Correspond Handoff Datetime (E2) / Correspond Handback DateTime (H2)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-17 18:54:31"
$wsDeDe.Range("H2").Value = "2016-03-17 18:54:52"
